$d = $word.ActiveDocument

# --- Paragraph: "During high school, I was a member of the FIRST Robotics ..." ---
# Replace the whole paragraph sentence with the revised wording.
$oldRobotics = "During high school, I was a member of the FIRST Robotics Competition Team 2928. I work with teammates from diverse backgrounds and cultures. During the team meetings, we often discuss the design and assembly process. My focus on the task and effectively communicating helped our team scores improve in world competitions. During the world competition, I used my Chinese language skills to help my teammates communicate with teams from China. From that, I believe that working with diverse groups in the future will bring unique perspectives. Other experiences in the robotics competition will also help me to be more respectful of others and develop leadership skills. I find those skills fulfilling to help others and it helps me build relationships."
$newRobotics = "During high school, I was a member of the FIRST Robotics Competition Team 2928. I work with teammates from diverse backgrounds and cultures. During the team meetings, my focus on the task and effectively communicating helped our team scores improve in world competitions. During the world competition, I used my Chinese language skills to help my teammates communicate with teams from China. From that, I will bring unique perspectives work with diverse groups in the future. Other experiences in the robotics competition will also help me to be more respectful of others and develop leadership skills. I find those skills fulfilling to help others and it helps me build relationships with other teammates."

$d.Content.Find.Execute($oldRobotics, $false, $false, $false, $false, $false, $true, 1, $false, $newRobotics, 2) | Out-Null

# Find the paragraph that now holds this text, drop the leading tab run and
# give it a first-line indent instead (matching the rest of the document's
# paragraphs).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.IndexOf("During high school, I was a member of the FIRST Robotics") -ge 0) {
        $r = $p.Range
        $firstChar = $d.Range($r.Start, $r.Start + 1)
        if ($firstChar.Text -eq "`t") {
            $firstChar.Text = ""
        }
        $p.Format.FirstLineIndent = 36
        break
    }
}

# --- Paragraph: "During the last winter brea. I surees deplove ..." ---
# Fix the typos and extend the paragraph with the new closing sentences.
$oldWinter = "During the last winter brea. I surees deplove my personal website. When I doing this project. I lerarn a lot from the project. I "
$newWinter = "During the last winter break. I success deployed my personal website. When I do this project, I learn a lot from the project. I build my website from scratch to basic done. As I build it, I start to understand how industry setting up their website. Doing this personal project, I strength my programming skill to another level. I start to switch my mind to thinking about how and what bug does my website will be remain that I need to be fixed."

$d.Content.Find.Execute($oldWinter, $false, $false, $false, $false, $false, $true, 1, $false, $newWinter, 2) | Out-Null
